$d = $word.ActiveDocument

# Locate the sentence to edit precisely via a duplicated range so we do not
# disturb $d.Content itself.
$target = $d.Content.Duplicate
$target.Find.Execute("la production des listings", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Turn on revision tracking just for this edit: replacing text while tracking
# changes is on forces Word to keep the surrounding / deleted / inserted text
# as separate runs instead of silently re-merging them into the original
# run. That mirrors how this paragraph ends up split into three <w:r> runs.
$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true

$target.Text = "les colis"

# Restore the original tracking setting and accept the two revisions we just
# created (the deletion of the old wording and the insertion of the new one)
# so the final document contains plain runs with no tracked-change markup,
# matching the committed result.
$d.TrackRevisions = $wasTracking
for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $d.Revisions.Item($i).Accept()
}
